# Add a new diary entry row to the end of the development diary table,
# describing the day's work (11/12/2021 - Stage 4 - 3 Hours).

$d = $word.ActiveDocument
$t = $d.Tables.Item(1)

# Add a new row at the end of the table. Word automatically inherits the
# cell formatting (shading, borders, widths, paragraph properties) from
# the preceding row.
$newRow = $t.Rows.Add()
$n = $t.Rows.Count

$t.Cell($n, 1).Range.Text = "11/12/2021"
$t.Cell($n, 2).Range.Text = "3 Hours"
$t.Cell($n, 3).Range.Text = "Stage 4"
$t.Cell($n, 4).Range.Text = "Changed the approach to arguments for queued functions by switching the system to a structure based system. This means that the values are not dropped as was the case with the char** implementation. However, a new issue has emerged with this implementation in which all queued actions seem to store the same data as the first queued action, likely due to some issue with the creation of new structs."
